# PEMstr.pptx edit script
# Commit message: "Update doc/slides with ignoreEOF"
#
# This script applies three logical changes described by the target diff:
#   1. Refresh the "datetimeFigureOut" date field placeholders from
#      4/10/13 -> 5/16/13 (Slide Master + every Slide Layout; the Notes
#      Master copy of this placeholder is not reachable/writable through
#      this COM surface, see note near the bottom).
#   2. Slide 23 ("Specifying Buffering"): re-wrap two bullets so each is a
#      single run instead of two runs split mid-sentence.
#   3. Slide 24 ("Some other useful options"): split the trailing
#      " (for named pipe creation)" run in two, and append a new bullet
#      paragraph describing "ignoreEOF (to allow multiple stream inputs)".

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Date placeholder fields: 4/10/13 -> 5/16/13
# ---------------------------------------------------------------------
$oldDate = "4/10/13"
$newDate = "5/16/13"

$master = $p.SlideMaster

# 1a. Slide Master itself
for ($j = 1; $j -le $master.Shapes.Count; $j++) {
    $sh = $master.Shapes.Item($j)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# 1b. Every Slide Layout under the master
$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $sh = $layout.Shapes.Item($j)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# 1c. Notes Master (present in the source deck, but this host's NotesMaster
#     shapes do not accept writes via COM - attempted here defensively in
#     case that changes; harmless no-op otherwise).
$notesMaster = $p.NotesMaster
for ($j = 1; $j -le $notesMaster.Shapes.Count; $j++) {
    $sh = $notesMaster.Shapes.Item($j)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# ---------------------------------------------------------------------
# 2. Slide 23 - "Specifying Buffering" content placeholder
# ---------------------------------------------------------------------
$slide23 = $p.Slides.Item(23)
$shape23 = $slide23.Shapes.Item(2)
$tr23 = $shape23.TextFrame.TextRange

# Paragraph 4: "Size of queue before we send back " + "pause" -> one run
$para4 = $tr23.Paragraphs(4, 1)
$para4Run1 = $para4.Runs(1, 1)
$para4Run1.Text = "Size of queue before we send back pause"
$para4Run2 = $para4.Runs(2, 1)
$para4Run2.Text = ""

# Paragraph 5: "# of buffers used for multi-buffering reads" + " " -> one run
$para5 = $tr23.Paragraphs(5, 1)
$para5Run1 = $para5.Runs(1, 1)
$para5Run1.Text = "# of buffers used for multi-buffering reads "
$para5Run2 = $para5.Runs(2, 1)
$para5Run2.Text = ""

# ---------------------------------------------------------------------
# 3. Slide 24 - "Some other useful options" content placeholder
# ---------------------------------------------------------------------
$slide24 = $p.Slides.Item(24)
$shape24 = $slide24.Shapes.Item(2)
$tr24 = $shape24.TextFrame.TextRange

# Last paragraph ("tmpDir (for named pipe creation)") currently has 2 runs;
# split the second run into " (for named pipe creation" + ")"
$lastPara = $tr24.Paragraphs(5, 1)
$lastParaRun2 = $lastPara.Runs(2, 1)
$lastParaRun2.Text = " (for named pipe creation"
$lastParaRun2.InsertAfter(")")

# Append a brand new bullet paragraph after it, built up run by run so the
# formatting boundaries match the target structure.
$tr24b = $shape24.TextFrame.TextRange
$tr24b.InsertAfter([char]13 + "ignoreEOF")

$tr24c = $shape24.TextFrame.TextRange
$newPara = $tr24c.Paragraphs(6, 1)
$newParaRun1 = $newPara.Runs(1, 1)
$newParaRun1.InsertAfter(" (to allow ")

$tr24d = $shape24.TextFrame.TextRange
$newPara2 = $tr24d.Paragraphs(6, 1)
$newParaRun2 = $newPara2.Runs(2, 1)
$newParaRun2.InsertAfter("multiple stream inputs")

$tr24e = $shape24.TextFrame.TextRange
$newPara3 = $tr24e.Paragraphs(6, 1)
$newParaRun3 = $newPara3.Runs(3, 1)
$newParaRun3.InsertAfter(")")

Write-Output "Edits applied"
